# Update clutter change order variables: flip several 1 -> 0 values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E5").Value = 0
$ws.Range("F12").Value = 0
$ws.Range("D13").Value = 0
$ws.Range("G19").Value = 0
$ws.Range("G21").Value = 0
$ws.Range("F31").Value = 0
$ws.Range("E33").Value = 0
$ws.Range("D46").Value = 0
$ws.Range("G49").Value = 0
$ws.Range("E52").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("F71").Value = 0
$ws.Range("G72").Value = 0

# Update the view to reflect scrolled position/selection
$ws.Range("K69").Select()
$excel.ActiveWindow.ScrollRow = 49
$excel.ActiveWindow.ScrollColumn = 1
